$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.651.17'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '''1.885.70'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("D4").Value = '''0.9997'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''246.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").Value = '''1.0000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").Value = '''0.4733'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '''0.2895'
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("D10").Value = '''22.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.24%  '
$ws.Range("D11").Value = '''100.01'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.56%  '
$ws.Range("D12").Value = '''0.7627'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.08%  '
$ws.Range("D13").Value = '''0.07825'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("D14").Value = '''1.881.81'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").Value = '''5.246'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").Value = '''284.96'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("D17").Value = '''30.602.20'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '''13.23'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").Value = '''0.000007534'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("D20").Value = '''0.9992'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("D21").Value = '''2.126.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.71%  '
$ws.Range("D22").Value = '''5.363'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.17%  '
$ws.Range("D23").Value = '''0.9997'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").Value = '''6.434'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.02%  '
$ws.Range("D25").Value = '''9.179'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.36%  '
$ws.Range("D26").Value = '''163.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("D27").Value = '''19.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.81%  '
$ws.Range("D28").Value = '''1.914'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").Value = '''0.09748'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  -1.05%  '
$ws.Range("D31").Value = '''1.502'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("D32").Value = '''4.257'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("D33").Value = '''4.192'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").Value = '''0.04855'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.91%  '
$ws.Range("D35").Value = '''1.132'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.49%  '
$ws.Range("D36").Value = '''0.7003'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.40%  '
$ws.Range("D37").Value = '''2.786'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.34%  '
$ws.Range("E38").Value = '  +0.77%  '
$ws.Range("D39").Value = '''2.882'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.76%  '
$ws.Range("D40").Value = '''6.326'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.64%  '
$ws.Range("D41").Value = '''75.53'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.77%  '
$ws.Range("D42").Value = '''1.978'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.26%  '
$ws.Range("D43").Value = '''0.4257'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("D44").Value = '''0.8413'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("D45").Value = '''0.9998'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''9.928'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.44%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '''101.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.46%  '
$ws.Range("D48").Value = '''7.040'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("D49").Value = '''35.34'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("D50").Value = '''0.05795'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").Value = '''0.3962'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.20%  '
